$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in category ("NoProject") for the no-project runs in rows 9-11
$ws.Range("E9").Value = "NoProject"
$ws.Range("E10").Value = "NoProject"
$ws.Range("E11").Value = "NoProject"

# Update the saved selection on the sheet to G24 (matches author's cursor position)
$ws.Activate()
$ws.Range("G24").Select()
